# SCD0020-001 - Admin SLN menambahkan data hari libur
# Update TC_ID Excel SCD0017 until SCD0025 and Update TC_ID Solution SCD0006 until SCD0025

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from SCD0307 -> SCD0020
$ws.Name = "SCD0020"

# Update the TC_ID cell (B2) from DGS-322 -> SCD0020-001
$ws.Range("B2").Value = "SCD0020-001"

# Column B is best-fit to its content; widen it to fit the longer TC_ID text
$ws.Columns.Item(2).ColumnWidth = 11.5

# Move/record the active selection to B3 (matches the saved view state)
$ws.Range("B3").Select() | Out-Null

# L2/N2 hold volatile TEXT(TODAY()+2, ...) formulas - they recompute
# automatically against the current date on recalculation/save, so no
# explicit write is needed here.
